# Split Database class into order, invoice, product, and user databases
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mis-cased status strings to the canonical "Done"
$ws.Range("G9").Value = "Done"
$ws.Range("G30").Value = "Done"

# Add assignee + completion status for the newly broken-out PANE sub-tasks
$ws.Range("F31").Value = "jon"
$ws.Range("G31").Value = "Done"

$ws.Range("F32").Value = "jon"
$ws.Range("G32").Value = "Done"

$ws.Range("F33").Value = "Alex"
$ws.Range("G33").Value = "Done"

$ws.Range("F34").Value = "Vinny"
$ws.Range("G34").Value = "Done"

# Update the view to reflect scrolling to the newly edited rows
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("F35").Select()
